{"js": "// Load all top-level body paragraphs (document has no tables, so body.paragraphs\n// covers everything we need to touch).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- 1) In-place text replacements (paragraph stays, text changes) -------\n// Keyed by the paragraph's 0-based index in the ORIGINAL document.\nconst replacements = {\n  3: \"Date: August 09, 2025\",\n  7: \"To help users manage their tasks efficiently.\",\n  16: \"A simple and intuitive to-do list app designed to help users manage tasks efficiently.\",\n  17: \"The app supports task creation, categorization, deadlines, and reminders.\",\n  20: \"3.1 Task Management [0]\",\n  22: \"Users can create, edit, and delete tasks.\",\n  26: \"- Implement task creation functionality\",\n  27: \"- Develop task editing features\",\n  28: \"- Design task deletion functionality\",\n  31: \"3.2 Task Scheduling [1]\",\n  33: \"Users can assign due dates and set reminders.\",\n  37: \"- Implement due date functionality for tasks.\",\n  38: \"- Allow users to set reminders for tasks with customizable notification times.\",\n  39: \"- Display due dates and reminders prominently in the task view.\",\n  42: \"3.3 Task Categorization [2]\",\n  44: \"Users can organize tasks using categories or tags.\",\n  48: \"- Implement a system for creating and managing categories.\",\n  49: \"- Allow users to assign multiple categories to a single task.\",\n  50: \"- Provide a clear visual representation of categories and their associated tasks.\",\n  51: \"- Enable filtering and searching tasks by category.\",\n  52: \"3.4 Deadline Notification [3]\",\n  54: \"The app should notify users of upcoming deadlines.\",\n  58: \"- Send push notifications for tasks nearing their deadline.\",\n  59: \"- Allow users to customize notification frequency and time.\",\n  60: \"- Implement email notifications as an alternative to push notifications.\",\n  63: \"3.5 Data Synchronization [4]\",\n  65: \"Tasks should sync across devices using cloud storage.\",\n  69: \"- Implement cloud storage integration (e.g., AWS S3, Firebase).\",\n  70: \"- Develop synchronization logic to handle task creation, updates, and deletions.\",\n  71: \"- Ensure secure data transfer between devices and the cloud.\",\n  72: \"- Implement conflict resolution strategies for concurrent modifications.\",\n  73: \"3.6 System Scalability [5]\",\n  75: \"The system should support at least 10,000 users simultaneously.\",\n  79: \"- Implement load balancing across multiple servers.\",\n  80: \"- Optimize database queries for scalability.\",\n  81: \"- Utilize caching mechanisms to reduce database load.\"\n};\n\nfor (const idxStr of Object.keys(replacements)) {\n  const idx = parseInt(idxStr, 10);\n  items[idx].insertText(replacements[idx], Word.InsertLocation.replace);\n}\n\n// --- 2) Paragraphs removed outright ---------------------------------------\n// - old para 18 (\"Users can mark tasks as complete ...\") merged away\n// - trimmed bullet lists under 3.1 / 3.2 / 3.4 / 3.6 (count reduced)\n// - the whole old 3.7-3.10 block (superseded by nothing; section 4 follows 3.6)\nconst deletions = [\n  18,\n  29, 30,\n  40, 41,\n  61, 62,\n  82\n];\nfor (let i = 83; i <= 123; i++) {\n  deletions.push(i);\n}\n\nfor (const idx of deletions) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) In-place text replacements (paragraph stays, only its text changes) ---\n# Paragraphs collection is 1-based; these numbers correspond to the ORIGINAL\n# (pre-edit) document layout and are applied before any deletions happen, so\n# the numbering is still valid here.\n$replacements = @{\n    4   = \"Date: August 09, 2025\"\n    8   = \"To help users manage their tasks efficiently.\"\n    17  = \"A simple and intuitive to-do list app designed to help users manage tasks efficiently.\"\n    18  = \"The app supports task creation, categorization, deadlines, and reminders.\"\n    21  = \"3.1 Task Management [0]\"\n    23  = \"Users can create, edit, and delete tasks.\"\n    27  = \"- Implement task creation functionality\"\n    28  = \"- Develop task editing features\"\n    29  = \"- Design task deletion functionality\"\n    32  = \"3.2 Task Scheduling [1]\"\n    34  = \"Users can assign due dates and set reminders.\"\n    38  = \"- Implement due date functionality for tasks.\"\n    39  = \"- Allow users to set reminders for tasks with customizable notification times.\"\n    40  = \"- Display due dates and reminders prominently in the task view.\"\n    43  = \"3.3 Task Categorization [2]\"\n    45  = \"Users can organize tasks using categories or tags.\"\n    49  = \"- Implement a system for creating and managing categories.\"\n    50  = \"- Allow users to assign multiple categories to a single task.\"\n    51  = \"- Provide a clear visual representation of categories and their associated tasks.\"\n    52  = \"- Enable filtering and searching tasks by category.\"\n    53  = \"3.4 Deadline Notification [3]\"\n    55  = \"The app should notify users of upcoming deadlines.\"\n    59  = \"- Send push notifications for tasks nearing their deadline.\"\n    60  = \"- Allow users to customize notification frequency and time.\"\n    61  = \"- Implement email notifications as an alternative to push notifications.\"\n    64  = \"3.5 Data Synchronization [4]\"\n    66  = \"Tasks should sync across devices using cloud storage.\"\n    70  = \"- Implement cloud storage integration (e.g., AWS S3, Firebase).\"\n    71  = \"- Develop synchronization logic to handle task creation, updates, and deletions.\"\n    72  = \"- Ensure secure data transfer between devices and the cloud.\"\n    73  = \"- Implement conflict resolution strategies for concurrent modifications.\"\n    74  = \"3.6 System Scalability [5]\"\n    76  = \"The system should support at least 10,000 users simultaneously.\"\n    80  = \"- Implement load balancing across multiple servers.\"\n    81  = \"- Optimize database queries for scalability.\"\n    82  = \"- Utilize caching mechanisms to reduce database load.\"\n}\n\nforeach ($idx in $replacements.Keys) {\n    $d.Paragraphs.Item($idx).Range.Text = $replacements[$idx]\n}\n\n# --- 2) Remove the trailing 3.6-subtask line plus the whole old 3.7 .. 3.10\n#        block in one shot (paragraphs 83..124 are contiguous). ---\n$startP = $d.Paragraphs.Item(83)\n$endP = $d.Paragraphs.Item(124)\n$blockRange = $d.Range($startP.Range.Start, $endP.Range.End)\n$blockRange.Delete()\n\n# --- 3) Remove the remaining scattered paragraphs, highest index first so\n#        earlier (lower) paragraph numbers stay valid while we work. ---\n$singleDeletions = @(63, 62, 42, 41, 31, 30, 19) | Sort-Object -Descending\nforeach ($idx in $singleDeletions) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
